# Updates coin Price (D) and Volume(1h) (E) text values on Sheet1,
# matching the "Updated symbol list" GitHub Actions commit.
#
# The source cells are stored as text (t="inlineStr"/shared string), not
# numbers (e.g. "0.1720" must keep its trailing zero and "0.34%" must stay
# a literal percent-string, not get parsed into 0.0034). Setting .Value on a
# General-formatted cell with a numeric-looking string auto-converts it to a
# real number, so for every target cell we temporarily switch its
# NumberFormat to Text ("@"), assign the literal string, then restore the
# "Normal" style so the on-disk cell style/format matches the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '291.95'
Set-TextValue 'E2' '0.34%'

Set-TextValue 'D3' '30.99'
Set-TextValue 'E3' '0.84%'

Set-TextValue 'D4' '4.955'
Set-TextValue 'E4' '1.31%'

Set-TextValue 'D5' '0.07452'
Set-TextValue 'E5' '2.63%'

Set-TextValue 'D6' '2.235'
Set-TextValue 'E6' '-4.97%'

Set-TextValue 'D7' '7.742'
Set-TextValue 'E7' '1.04%'

Set-TextValue 'D8' '0.9197'
Set-TextValue 'E8' '2.49%'

Set-TextValue 'D9' '0.09286'
Set-TextValue 'E9' '16.81%'

Set-TextValue 'D10' '0.1720'
Set-TextValue 'E10' '3.01%'

Set-TextValue 'D11' '0.08326'
Set-TextValue 'E11' '1.69%'

Set-TextValue 'D12' '0.03229'
Set-TextValue 'E12' '4.75%'

Set-TextValue 'D13' '0.09932'
Set-TextValue 'E13' '-0.93%'

Set-TextValue 'D14' '0.001493'
Set-TextValue 'E14' '-0.22%'

Set-TextValue 'D15' '0.005742'
Set-TextValue 'E15' '-1.48%'

Set-TextValue 'D16' '3.477'
Set-TextValue 'E16' '0.07%'

Set-TextValue 'D17' '3.758'
Set-TextValue 'E17' '1.57%'

Set-TextValue 'D18' '2.129'
Set-TextValue 'E18' '2.47%'

Set-TextValue 'D19' '0.3332'
Set-TextValue 'E19' '0.40%'

Set-TextValue 'D20' '0.1302'
Set-TextValue 'E20' '0.32%'

Set-TextValue 'D21' '4.176'
Set-TextValue 'E21' '5.21%'

Set-TextValue 'E22' '-8.06%'

Set-TextValue 'D23' '0.04502'
Set-TextValue 'E23' '-0.45%'

Set-TextValue 'D24' '0.001216'
Set-TextValue 'E24' '0.49%'

Set-TextValue 'D25' '0.004261'
Set-TextValue 'E25' '-3.46%'

Set-TextValue 'D26' '0.0001298'

Set-TextValue 'D27' '0.0003385'
Set-TextValue 'E27' '-0.19%'

Set-TextValue 'D39' '0.01623'
Set-TextValue 'E39' '2.17%'

Set-TextValue 'D40' '0.04578'
Set-TextValue 'E40' '4.60%'

Set-TextValue 'D41' '0.007407'
Set-TextValue 'E41' '1.17%'

Set-TextValue 'D42' '0.009817'
Set-TextValue 'E42' '-1.96%'

Set-TextValue 'D43' '0.1357'

Set-TextValue 'D44' '0.002216'
Set-TextValue 'E44' '9.84%'

Set-TextValue 'D45' '0.009637'
Set-TextValue 'E45' '1.37%'

Set-TextValue 'D46' '0.00006091'
Set-TextValue 'E46' '6.50%'

Set-TextValue 'E47' '-0.35%'

Set-TextValue 'D48' '2.655'
Set-TextValue 'E48' '18.39%'

Set-TextValue 'D49' '0.001995'
Set-TextValue 'E49' '-31.14%'

Set-TextValue 'E50' '-0.35%'

Set-TextValue 'D51' '0.0001995'
Set-TextValue 'E51' '-0.35%'
